$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 158, shifting existing rows 158:173 down to 159:174
$ws.Rows.Item(158).Insert()

# Populate the new row 158 with the new record
$ws.Cells.Item(158, 1).Value = 10
$ws.Cells.Item(158, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(158, 3).Value = "La Araucanía"
$ws.Cells.Item(158, 4).Value = 44748
$ws.Cells.Item(158, 5).Value = 9
$ws.Cells.Item(158, 6).Value = 100114007
$ws.Cells.Item(158, 7).Value = "Jengibre"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 20
$ws.Cells.Item(158, 11).Value = 20000
$ws.Cells.Item(158, 12).Value = 20000
$ws.Cells.Item(158, 13).Value = 20000
$ws.Cells.Item(158, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(158, 15).Value = "Perú"
$ws.Cells.Item(158, 16).Value = 1538
$ws.Cells.Item(158, 17).Value = 13
$ws.Cells.Item(158, 18).Value = "Hortaliza"
